$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add Tags values for rows 5-8
$ws.Range("G5").Value = "1,1"
$ws.Range("G6").Value = "1,2"
$ws.Range("G7").Value = "1,3"
$ws.Range("G8").Value = "1,4"

# Update the selection shown in the sheet view
$ws.Activate()
$ws.Range("A5:J8").Select()
